# Generate Report for Handoff
# Update the localization status report: the "b.md" file has been
# re-handed-off (a new handoff xliff was generated), so its status moves
# from "Handed back: in sync with en-US" to "Ready for handoff" on every
# sheet, along with the new handoff file name/datetime and a warning
# that the previous handback is now stale.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc8252dbce87c16a6647be7c64c8c789277ebfb7/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/252c1caba148fb6b52238168588b3375094d3923/e2e/b.md."

# ----- Overview sheet -----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus
$ov.Range("G3").Value = "2016-08-29 16:41:31"

# ----- zh-cn sheet (row 3 = b.md) -----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $newStatus
# Leading apostrophe forces this to be stored as text "False" instead of
# being auto-coerced into a native boolean value.
$zh.Range("F3").Value = "'False"
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-29 16:41:26"
$zh.Range("P3").Value = $errorDetail
$zh.Columns.Item(16).ColumnWidth = 40

# ----- de-de sheet (row 3 = b.md) -----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $newStatus
$de.Range("F3").Value = "'False"
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = "2016-08-29 16:41:31"
$de.Range("P3").Value = $errorDetail
$de.Columns.Item(16).ColumnWidth = 40
